$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A, shifting the existing query/dbExcel/WebExcel columns
# (old A:D) one place to the right (new B:E). Values, styles and column
# widths of the old columns are carried along automatically by Insert().
$ws.Columns("A:A").Insert()

# New first column: tab name / label for the query that now lives in
# columns B/C.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Replace the old "Trials" query (now in column B, row 2) with the new
# "Cases" query text.
$casesQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "BLACK_OR_AFRICAN_AMERICAN"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@
$ws.Range("B2").Value = $casesQuery

# Replace the old stat query (now in column C, row 2) with the new stat
# query text.
$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "BLACK_OR_AFRICAN_AMERICAN"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@
$ws.Range("C2").Value = $statQuery

# New narrow first column (~8.8 chars, auto-fit to the "CasesTab"/"TabName"
# labels).
$ws.Columns("A:A").ColumnWidth = 7.92

# Row 2 grows taller to fit the longer, multi-line queries.
$ws.Rows("2:2").RowHeight = 174

# Selection moves down to B5, matching the saved view state.
$ws.Range("B5").Select()
